# Insert two new weekly price rows at the top of the Pera (pear) data block,
# pushing the existing rows 1063..1134 down to 1065..1136.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows right before the current row 1063.
$ws.Range("A1063:T1064").EntireRow.Insert()

# --- New row 1063: Packham's Triumph / Calibre 80 / Provincia de Curico ---
$ws.Cells.Item(1063, 1).Value = 9
$ws.Cells.Item(1063, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1063, 3).Value = "Metropolitana"
$ws.Cells.Item(1063, 4).Value = 45106
$ws.Cells.Item(1063, 5).Value = 13
$ws.Cells.Item(1063, 6).Value = "Fruta"
$ws.Cells.Item(1063, 7).Value = 100104
$ws.Cells.Item(1063, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1063, 9).Value = 100104005
$ws.Cells.Item(1063, 10).Value = "Pera"
$ws.Cells.Item(1063, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1063, 12).Value = "Calibre 80"
$ws.Cells.Item(1063, 13).Value = 350
$ws.Cells.Item(1063, 14).Value = 17000
$ws.Cells.Item(1063, 15).Value = 17000
$ws.Cells.Item(1063, 16).Value = 17000
$ws.Cells.Item(1063, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(1063, 18).Value = "Provincia de Curic" + [char]243
$ws.Cells.Item(1063, 19).Value = 944
$ws.Cells.Item(1063, 20).Value = 18

# --- New row 1064: Packham's Triumph / Primera / Region de O'Higgins ---
$ws.Cells.Item(1064, 1).Value = 9
$ws.Cells.Item(1064, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1064, 3).Value = "Metropolitana"
$ws.Cells.Item(1064, 4).Value = 45106
$ws.Cells.Item(1064, 5).Value = 13
$ws.Cells.Item(1064, 6).Value = "Fruta"
$ws.Cells.Item(1064, 7).Value = 100104
$ws.Cells.Item(1064, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1064, 9).Value = 100104005
$ws.Cells.Item(1064, 10).Value = "Pera"
$ws.Cells.Item(1064, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1064, 12).Value = "Primera"
$ws.Cells.Item(1064, 13).Value = 590
$ws.Cells.Item(1064, 14).Value = 13000
$ws.Cells.Item(1064, 15).Value = 14000
$ws.Cells.Item(1064, 16).Value = 13475
$ws.Cells.Item(1064, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(1064, 18).Value = "Regi" + [char]243 + "n de O'Higgins"
$ws.Cells.Item(1064, 19).Value = 749
$ws.Cells.Item(1064, 20).Value = 18
